$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 14.26366496892708
$ws.Range("C2").Value = 11.10383122448296
$ws.Range("D2").Value = 6.826510824745297
$ws.Range("E2").Value = 12.73965231580022
$ws.Range("F2").Value = 44.35775937862943
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 29.52759519459227
$ws.Range("J2").Value = 10.37009018898844
$ws.Range("K2").Value = 15.3236777357165
$ws.Range("N2").Value = 21.97709096228471
$ws.Range("B3").Value = 14.06300550304561
$ws.Range("C3").Value = 10.95478535355923
$ws.Range("D3").Value = 6.799054472679478
$ws.Range("E3").Value = 12.69667285973865
$ws.Range("F3").Value = 44.30839309715828
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 29.56263075875016
$ws.Range("J3").Value = 10.37806476201765
$ws.Range("K3").Value = 15.19114164693403
$ws.Range("N3").Value = 22.034024535239
$ws.Range("B4").Value = 13.94228964832539
$ws.Range("C4").Value = 10.86551881285956
$ws.Range("D4").Value = 6.783485936586973
$ws.Range("E4").Value = 12.67304878075003
$ws.Range("F4").Value = 44.28792935998083
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 29.5895240158062
$ws.Range("J4").Value = 10.38465996792802
$ws.Range("K4").Value = 15.1130951457311
$ws.Range("N4").Value = 22.0709082416699
$ws.Range("B5").Value = 13.8937867910274
$ws.Range("C5").Value = 10.82975257494245
$ws.Range("D5").Value = 6.77747097645235
$ws.Range("E5").Value = 12.66412392969218
$ws.Range("F5").Value = 44.28206962835082
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 29.60183383996501
$ws.Range("J5").Value = 10.38777472127817
$ws.Range("K5").Value = 15.08216023924443
$ws.Range("N5").Value = 22.08642363047197
$ws.Range("B6").Value = 13.88577643082104
$ws.Range("C6").Value = 10.82385179452109
$ws.Range("D6").Value = 6.77649222769657
$ws.Range("E6").Value = 12.66268455484101
$ws.Range("F6").Value = 44.28124643240019
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 29.60395936863678
$ws.Range("J6").Value = 10.38831771997561
$ws.Range("K6").Value = 15.07707694989734
$ws.Range("N6").Value = 22.08902925854161
$ws.Range("B7").Value = 13.9416326454539
$ws.Range("C7").Value = 10.86503392572703
$ws.Range("D7").Value = 6.783403476872047
$ws.Range("E7").Value = 12.67292556576712
$ws.Range("F7").Value = 44.28784029180669
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 29.58968456557199
$ws.Range("J7").Value = 10.38470024516441
$ws.Range("K7").Value = 15.11267438447533
$ws.Range("N7").Value = 22.07111552320043
$ws.Range("B8").Value = 14.19399939274679
$ws.Range("C8").Value = 11.05200111653298
$ws.Range("D8").Value = 6.816779307126867
$ws.Range("E8").Value = 12.72426331083814
$ws.Range("F8").Value = 44.33869651950432
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 29.53855712255308
$ws.Range("J8").Value = 10.3724872563529
$ws.Range("K8").Value = 15.27730762393872
$ws.Range("N8").Value = 21.99632222618136
$ws.Range("B9").Value = 14.70580014342942
$ws.Range("C9").Value = 11.43444057971924
$ws.Range("D9").Value = 6.892225652731645
$ws.Range("E9").Value = 12.84654757021663
$ws.Range("F9").Value = 44.51634485075301
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 29.48109831075305
$ws.Range("J9").Value = 10.36201449057142
$ws.Range("K9").Value = 15.62510029526656
$ws.Range("N9").Value = 21.8649087720981
$ws.Range("B10").Value = 15.08829679470492
$ws.Range("C10").Value = 11.72227335233791
$ws.Range("D10").Value = 6.953409199125095
$ws.Range("E10").Value = 12.94907084604606
$ws.Range("F10").Value = 44.69396872344528
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 29.46510269913158
$ws.Range("J10").Value = 10.36252830157753
$ws.Range("K10").Value = 15.89375124522732
$ws.Range("N10").Value = 21.77761988052275
$ws.Range("B11").Value = 15.26290588042936
$ws.Range("C11").Value = 11.85411381437447
$ws.Range("D11").Value = 6.982413004548213
$ws.Range("E11").Value = 12.99834462361385
$ws.Range("F11").Value = 44.7848799742106
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 29.46353892312525
$ws.Range("J11").Value = 10.36454038056826
$ws.Range("K11").Value = 16.01835470078394
$ws.Range("N11").Value = 21.73991257929061
$ws.Range("B12").Value = 15.32904780528343
$ws.Range("C12").Value = 11.90411971229129
$ws.Range("D12").Value = 6.993557243679178
$ws.Range("E12").Value = 13.01737107500141
$ws.Range("F12").Value = 44.82074577326814
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 29.46376911143742
$ws.Range("J12").Value = 10.36555744079665
$ws.Range("K12").Value = 16.06584309148684
$ws.Range("N12").Value = 21.72592099939134
$ws.Range("B13").Value = 15.31480305093255
$ws.Range("C13").Value = 11.89334722738051
$ws.Range("D13").Value = 6.991150084711481
$ws.Range("E13").Value = 13.01325723026889
$ws.Range("F13").Value = 44.81295764432826
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 29.46368295213011
$ws.Range("J13").Value = 10.36532706256098
$ws.Range("K13").Value = 16.05560275127095
$ws.Range("N13").Value = 21.72892156171937
$ws.Range("B14").Value = 15.26834736432934
$ws.Range("C14").Value = 11.85822648901668
$ws.Range("D14").Value = 6.983326661140718
$ws.Range("E14").Value = 12.99990264279201
$ws.Range("F14").Value = 44.78780190859428
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 29.46354137605006
$ws.Range("J14").Value = 10.3646189454116
$ws.Range("K14").Value = 16.02225574314938
$ws.Range("N14").Value = 21.738755726535
$ws.Range("B15").Value = 15.23989274279142
$ws.Range("C15").Value = 11.83672309905011
$ws.Range("D15").Value = 6.978555351616687
$ws.Range("E15").Value = 12.99177008735628
$ws.Range("F15").Value = 44.77258033548288
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 29.46356177069913
$ws.Range("J15").Value = 10.36421840877166
$ws.Range("K15").Value = 16.0018680961833
$ws.Range("N15").Value = 21.74481684837894
$ws.Range("B16").Value = 15.07689222998712
$ws.Range("C16").Value = 11.71367122875339
$ws.Range("D16").Value = 6.951536681747007
$ws.Range("E16").Value = 12.94590275747106
$ws.Range("F16").Value = 44.68822972769236
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 29.46531993403223
$ws.Range("J16").Value = 10.36243254447877
$ws.Range("K16").Value = 15.88565282609656
$ws.Range("N16").Value = 21.78012437199725
$ws.Range("B17").Value = 14.97700502480317
$ws.Range("C17").Value = 11.6383790985823
$ws.Range("D17").Value = 6.935256456710169
$ws.Range("E17").Value = 12.91843169925913
$ws.Range("F17").Value = 44.63906320390229
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 29.46786242424789
$ws.Range("J17").Value = 10.36179206946068
$ws.Range("K17").Value = 15.81494279404481
$ws.Range("N17").Value = 21.80229659451923
$ws.Range("B18").Value = 14.91961378587487
$ws.Range("C18").Value = 11.59516102893571
$ws.Range("D18").Value = 6.926003184418566
$ws.Range("E18").Value = 12.90287985714585
$ws.Range("F18").Value = 44.61173640912175
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 29.46986247747225
$ws.Range("J18").Value = 10.36159111503046
$ws.Range("K18").Value = 15.77449991463594
$ws.Range("N18").Value = 21.81523782034982
$ws.Range("B19").Value = 14.90019466691235
$ws.Range("C19").Value = 11.58054475710671
$ws.Range("D19").Value = 6.922889412683077
$ws.Range("E19").Value = 12.89765733750419
$ws.Range("F19").Value = 44.60264801366636
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 29.47063197493974
$ws.Range("J19").Value = 10.36155184753865
$ws.Range("K19").Value = 15.76084697477173
$ws.Range("N19").Value = 21.81965185784631
$ws.Range("B20").Value = 14.98763233659749
$ws.Range("C20").Value = 11.64638532908087
$ws.Range("D20").Value = 6.936978111127817
$ws.Range("E20").Value = 12.92133037292787
$ws.Range("F20").Value = 44.64419858076985
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 29.46753611623717
$ws.Range("J20").Value = 10.36184292464851
$ws.Range("K20").Value = 15.8224467267152
$ws.Range("N20").Value = 21.79991683273964
$ws.Range("B21").Value = 15.28199246895723
$ws.Range("C21").Value = 11.86854049603316
$ws.Range("D21").Value = 6.985620276562568
$ws.Range("E21").Value = 13.00381532801872
$ws.Range("F21").Value = 44.79515180629805
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 29.46356063690285
$ws.Range("J21").Value = 10.36482001781674
$ws.Range("K21").Value = 16.03204264617879
$ws.Range("N21").Value = 21.73585939794235
$ws.Range("B22").Value = 15.47446825684587
$ws.Range("C22").Value = 12.01418052415597
$ws.Range("D22").Value = 7.01834637222151
$ws.Range("E22").Value = 13.05986060285609
$ws.Range("F22").Value = 44.90219223332225
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 29.4657558699903
$ws.Range("J22").Value = 10.36825244356229
$ws.Range("K22").Value = 16.17077745272061
$ws.Range("N22").Value = 21.69566894449103
$ws.Range("B23").Value = 15.3717531195346
$ws.Range("C23").Value = 11.9364246142187
$ws.Range("D23").Value = 7.000796696044138
$ws.Range("E23").Value = 13.02975661866379
$ws.Range("F23").Value = 44.84430073418738
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 29.46414544946424
$ws.Range("J23").Value = 10.36628468784233
$ws.Range("K23").Value = 16.09658516160235
$ws.Range("N23").Value = 21.71696621708448
$ws.Range("B24").Value = 14.98282761402672
$ws.Range("C24").Value = 11.64276549551
$ws.Range("D24").Value = 6.936199419063296
$ws.Range("E24").Value = 12.9200191290556
$ws.Range("F24").Value = 44.64187394859015
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 29.46768196312948
$ws.Range("J24").Value = 10.3618194119634
$ws.Range("K24").Value = 15.8190535438252
$ws.Range("N24").Value = 21.80099211825748
$ws.Range("B25").Value = 14.56593065830292
$ws.Range("C25").Value = 11.32957514812381
$ws.Range("D25").Value = 6.870779082263987
$ws.Range("E25").Value = 12.81120082028426
$ws.Range("F25").Value = 44.45997784860892
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 29.49204571029554
$ws.Range("J25").Value = 10.36340517480068
$ws.Range("K25").Value = 15.52855282792173
$ws.Range("N25").Value = 21.89883054336965
